$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet (column G holds "Recorded By" values).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ($text -ne $null -and $text.StartsWith("System, ")) {
        $parts = $text.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $reversed = $parts[($parts.Length - 1)..0]
        $newText = [string]::Join(", ", $reversed)
        $cell.Value = $newText
    }
}
